# Sprint 4 burndown chart update: US#5 row removed, new daily-progress
# values entered for the remaining four user stories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (US#1 Page:Eliott) : columns J..AD (10..30) ---
$row2 = @(1,0,0,0,1,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 10 + $i).Value = $row2[$i]
}

# --- Row 3 (US#2 Navigation:Eliott+Max) : J stays "0.5", K..AD -> 0 ---
$ws.Cells.Item(3, 10).Text = "0.5"
for ($col = 11; $col -le 30; $col++) {
    $ws.Cells.Item(3, $col).Value = 0
}

# --- Row 4 (US#3 Class 'User':Vincenzo) : J stays 4, K..AD -> 0 ---
$ws.Cells.Item(4, 10).Value = 4
for ($col = 11; $col -le 30; $col++) {
    $ws.Cells.Item(4, $col).Value = 0
}

# --- Row 5 (US#4 SQL:Mika+Karim+Al) : J stays 4, K..AD -> new values ---
$row5 = @(4,0,0,0,2,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 10 + $i).Value = $row5[$i]
}

# --- Row 6 (US#5 API:Al) : the whole user story row is removed ---
for ($col = 1; $col -le 9; $col++) {
    $ws.Cells.Item(6, $col).ClearContents()
}

# Refresh the burndown chart so its cached series reflect the new data.
$ws.ChartObjects(1).Chart.Refresh()
